$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C6").Value = "https://ipfs.infura.io:5001/api/v0/block/get?arg=QmQv9RDA4LPLsQzFUncHmtw8kWWJBrnqyxQm53marjXPFY"
